$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.269.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "'3.544.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'615.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.91%  "
$ws.Range("D6").Value = "'186.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "'0.659"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "'53.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("D13").Value = "'9.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "'4.096.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "'615.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.10%  "
$ws.Range("D16").Value = "'70.163.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "'12.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "'19.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").Value = "'3.546.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").Value = "'17.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "'104.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.54%  "
$ws.Range("D24").Value = "'4.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").Value = "'5.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").Value = "'3.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.11%  "
$ws.Range("D27").Value = "'11.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").Value = "'9.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.49%  "
$ws.Range("D29").Value = "'33.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.37%  "
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("D31").Value = "'12.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "'64.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "'3.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.15%  "
$ws.Range("E35").Value = "  -3.08%  "
$ws.Range("D36").Value = "'533.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.42%  "
$ws.Range("D37").Value = "'0.401"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'37.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("D40").Value = "'3.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.77%  "
$ws.Range("D41").Value = "'0.0₃0780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.40%  "
$ws.Range("D42").Value = "'3.525.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  +4.62%  "
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  +3.75%  "
$ws.Range("D47").Value = "'3.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  -2.63%  "
$ws.Range("D51").Value = "'134.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
